# Auto-generated script applying the cryptos.xlsx price/volume update
# (commit: "Updated cryptos list on Thu Jun  6 21:29:50 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cells (coin names, links, percentages, and price strings that
# already contain two '.' separators so Excel can't coerce them to numbers)
# can be assigned directly.
$ws.Range('D2').Value = '70.710.73'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').Value = '3.798.06'
$ws.Range('E3').Value = '  -1.44%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('E5').Value = '  +1.27%  '
$ws.Range('E6').Value = '  -1.43%  '
$ws.Range('D7').Value = '3.796.96'
$ws.Range('E7').Value = '  -1.44%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -1.15%  '
$ws.Range('E10').Value = '  -1.98%  '
$ws.Range('E11').Value = '  +0.60%  '
$ws.Range('E12').Value = '  -1.33%  '
$ws.Range('E13').Value = '  -2.08%  '
$ws.Range('E14').Value = '  -1.33%  '
$ws.Range('D15').Value = '4.440.26'
$ws.Range('E15').Value = '  -1.42%  '
$ws.Range('D16').Value = '3.781.26'
$ws.Range('E16').Value = '  -2.05%  '
$ws.Range('D17').Value = '70.740.60'
$ws.Range('E17').Value = '  -0.62%  '
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('E19').Value = '  -1.91%  '
$ws.Range('E20').Value = '  -2.08%  '
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('E22').Value = '  -4.79%  '
$ws.Range('E23').Value = '  +0.60%  '
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('E25').Value = '  -1.23%  '
$ws.Range('E26').Value = '  -1.92%  '
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('D28').Value = '3.951.00'
$ws.Range('E28').Value = '  -1.54%  '
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('E30').Value = '  -4.76%  '
$ws.Range('E31').Value = '  -2.76%  '
$ws.Range('E32').Value = '  -3.76%  '
$ws.Range('E33').Value = '  -4.17%  '
$ws.Range('E34').Value = '  -2.34%  '
$ws.Range('E35').Value = '  -2.43%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('B37').Value = 'RenzoRestakedETH'
$ws.Range('C37').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D37').Value = '3.769.48'
$ws.Range('E37').Value = '  -0.97%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('E38').Value = '  -2.49%  '
$ws.Range('E39').Value = '  -3.86%  '
$ws.Range('E40').Value = '  +1.71%  '
$ws.Range('E41').Value = '  -3.28%  '
$ws.Range('E42').Value = '  -2.10%  '
$ws.Range('E43').Value = '  -3.92%  '
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('E46').Value = '  +4.20%  '
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('E48').Value = '  +1.54%  '
$ws.Range('E49').Value = '  +0.16%  '
$ws.Range('E50').Value = '  -1.16%  '
$ws.Range('E51').Value = '  -1.85%  '

# Price cells whose text looks like a plain number (single '.') would be
# auto-coerced to a numeric type by the Value setter. Force the cell to
# Text format first so the literal string is preserved, then clear the
# number-format override so the cell's style stays at its original default.
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '704.66'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.53'
$ws.Range('D6').ClearFormats()
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.520'
$ws.Range('D9').ClearFormats()
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.40'
$ws.Range('D11').ClearFormats()
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.456'
$ws.Range('D12').ClearFormats()
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000253'
$ws.Range('D13').ClearFormats()
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.11'
$ws.Range('D14').ClearFormats()
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.13'
$ws.Range('D19').ClearFormats()
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.38'
$ws.Range('D20').ClearFormats()
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '495.50'
$ws.Range('D21').ClearFormats()
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.60'
$ws.Range('D22').ClearFormats()
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.728'
$ws.Range('D23').ClearFormats()
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.01'
$ws.Range('D24').ClearFormats()
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000144'
$ws.Range('D25').ClearFormats()
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.44'
$ws.Range('D27').ClearFormats()
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.33'
$ws.Range('D32').ClearFormats()
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.21'
$ws.Range('D33').ClearFormats()
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.07'
$ws.Range('D34').ClearFormats()
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').ClearFormats()
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '9.07'
$ws.Range('D38').ClearFormats()
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.93'
$ws.Range('D42').ClearFormats()
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.28'
$ws.Range('D43').ClearFormats()
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.000319'
$ws.Range('D46').ClearFormats()
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '164.68'
$ws.Range('D47').ClearFormats()
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '425.20'
$ws.Range('D48').ClearFormats()
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '48.77'
$ws.Range('D49').ClearFormats()
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.58'
$ws.Range('D50').ClearFormats()
